$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "IT Specialist (Remote)"
$ws.Range("E2").Value = "`$50,000 - `$60,000 a year"
$ws.Range("F2").Value = "Compensation Range: `$50,000 - `$60,000 per year dependent upon experience"

# Row 3
$ws.Range("A3").Value = "IT PROFESSIONAL 2"
$ws.Range("C3").Value = "Not Disclosed"
$ws.Range("E3").Value = "`$73,310 - `$109,641 a year"
$ws.Range("F3").Value = "No Description Available"

# Row 4
$ws.Range("A4").Value = "Supervisor, IT Service Desk"
$ws.Range("C4").Value = "Not Disclosed"
$ws.Range("E4").Value = "`$80,000 - `$90,000 a year"
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("A5").Value = "IT Help Desk Technician"
$ws.Range("C5").Value = "Manage user accounts. | Train new hires on infrastructure and operating system use. | Train users when systems or applications are upgraded. | Schedule requests to add, move or change hardware, software or phones. | Document procedures and maintain a library of information so that repeat questions or issues can be solved quickly"
$ws.Range("E5").Value = "`$30.08 - `$38.92 an hour"

# Row 6
$ws.Range("A6").Value = "IT Service Desk Specialist"
$ws.Range("C6").Value = "Not Disclosed"

# Row 7
$ws.Range("A7").Value = "IT Coordinator"
$ws.Range("C7").Value = "Proficient with Google products and Android (i.e. Google Business Suite) | Proficient with Apple products and iOS | Proficient with Windows products and operating software | Experience with HTML, Wordpress, Squarespace, and Shopify | Highly developed hardware/software troubleshooting techniques | Ability to provide L2/L3 technical support | Willingness to work a flexible schedule, including nights and weekends, based on the needs of the business | Excellent communication and problem solving skills"
$ws.Range("E7").Value = "`$55,000 - `$65,000 a year"
$ws.Range("F7").Value = "Marketing & Accounting Support"

# Row 8
$ws.Range("A8").Value = "IT Systems Administrator II – Denver International Airport"
$ws.Range("C8").Value = "FBI Background Check: FBI criminal background check is required for all positions at Denver International Airport (DEN). Employees are also required to report any felony convictions and/or moving violations to maintain this clearance and be eligible for continued employment. By position, a pre-employment physical/drug test may be required. | Snow/Emergency Duties: Denver International Airport is a 24/7/365 team operation. If weather conditions warrant or an emergency crisis occurs, all DEN employees can be required to work extended hours and/or shifts."
$ws.Range("E8").Value = "`$70,765 - `$116,762 a year"
$ws.Range("F8").Value = "No Description Available"

# Row 9
$ws.Range("A9").Value = "Service Desk Technician - IT"
$ws.Range("E9").Value = "From `$24.78 an hour"
$ws.Range("F9").Value = "Starting Hourly Rate: 24.78"

# Row 10
$ws.Range("A10").Value = "IT TECHNICIAN 4"
$ws.Range("E10").Value = "`$56,689.00 - `$83,666.16 a year"
$ws.Range("F10").Value = "No Description Available"

# Row 11
$ws.Range("A11").Value = "IT Support Specialist"
$ws.Range("C11").Value = "Not Disclosed"
$ws.Range("E11").Value = "Full-time"
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("A12").Value = "IT Helpdesk Specialist"
$ws.Range("E12").Value = "No Location Available"

# Row 13
$ws.Range("A13").Value = "Computer Systems Technician, Walt Disney World"
$ws.Range("E13").Value = "`$37.43 an hour"

# Row 14
$ws.Range("A14").Value = "IT Specialist"
$ws.Range("C14").Value = "Not Disclosed"
$ws.Range("E14").Value = "Full-time"
$ws.Range("F14").Value = "No Description Available"

# Row 15
$ws.Range("A15").Value = "Services Analyst - IT"
$ws.Range("C15").Value = "Annual salary will vary based on a candidate's skills, qualifications, experience, and other factors: `$84,000-`$126,000 | Annual bonus and incentive pay up to 10% | 401(k) match and annual company contribution | Medical, Dental and Vision Insurance | Life and disability insurance | Generous paid time off, including vacation, floating and fixed holidays and sick time | Maternity leave as well as paid bonding/primary caregiver leave or parental leave for the birth or adoption of a child or to care for an ill family member, as applicable (eligibility based on position) | Long Term Incentive Plan for eligible positions | Wellbeing programs such as tuition reimbursement, adoption assistance and fitness reimbursement | Referral bonus program | And much more"
$ws.Range("E15").Value = "`$84,000 - `$126,000 a year"
$ws.Range("F15").Value = "No Description Available"

# Row 16
$ws.Range("A16").Value = "IT Support Technician 1"
$ws.Range("E16").Value = "Full-time"
